$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FAT")

# Row 21: D21 200 -> "NA"
$ws.Range("D21").Value = "NA"

# Rows where D:J columns all become 0 (previously "NA" text in E:J, and a number in D)
$zeroRows = @(41, 43, 45, 46, 47, 49, 52, 54, 57, 59, 60, 62, 66, 72, 76)
foreach ($r in $zeroRows) {
    $ws.Range("D$r`:J$r").Value = 0
}

# Row 61: only D61 changes to 0 (E61:J61 already 0)
$ws.Range("D61").Value = 0

# Rows where D column becomes "NA"
$naRows = @(83, 89, 94, 100, 102)
foreach ($r in $naRows) {
    $ws.Range("D$r").Value = "NA"
}

# Row 101: J101 becomes "NA"
$ws.Range("J101").Value = "NA"
